# Correct Avionics cost per kW for FG, plus the related navigation/selection
# state that Excel records as part of the same save.

$wb = $excel.ActiveWorkbook

# --- kite sheet: correct the Avionics cost (avio.C, B8) and move the cursor ---
$kite = $wb.Worksheets.Item("kite")
$kite.Activate() | Out-Null
$kite.Range("B25").Select() | Out-Null
$kite.Range("B8").Value = 150000

# --- gStation sheet: correct a value (B1) and leave it as the active sheet ---
$gStation = $wb.Worksheets.Item("gStation")
$gStation.Activate() | Out-Null
$gStation.Range("E9").Select() | Out-Null
$gStation.Range("B1").Value = 100
